$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price and Volume columns to remain text so numeric-looking values
# (e.g. "1.00", "380.65") are not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "51.509.38"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").Value = "2.997.04"
$ws.Range("E3").Value = "  +2.67%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "380.65"
$ws.Range("E5").Value = "  +1.61%  "
$ws.Range("D6").Value = "103.32"
$ws.Range("E6").Value = "  +4.11%  "
$ws.Range("D7").Value = "0.545"
$ws.Range("E7").Value = "  +1.86%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +3.20%  "
$ws.Range("D10").Value = "36.70"
$ws.Range("E10").Value = "  +2.81%  "
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "0.0857"
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "3.461.24"
$ws.Range("E13").Value = "  +2.41%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "18.53"
$ws.Range("E14").Value = "  +3.16%  "
$ws.Range("E15").Value = "  +3.18%  "
$ws.Range("D16").Value = "2.991.50"
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("D17").Value = "11.31"
$ws.Range("E17").Value = "  -6.12%  "
$ws.Range("D18").Value = "0.989"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").Value = "51.533.84"
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("D20").Value = "3.09"
$ws.Range("E20").Value = "  +2.62%  "
$ws.Range("D21").Value = "12.51"
$ws.Range("E21").Value = "  +2.78%  "
$ws.Range("D22").Value = "0.0₃0963"
$ws.Range("E22").Value = "  +2.00%  "
$ws.Range("D23").Value = "70.25"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("D24").Value = "268.06"
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("D25").Value = "3.22"
$ws.Range("E25").Value = "  +2.45%  "
$ws.Range("D26").Value = "8.28"
$ws.Range("E26").Value = "  +6.23%  "
$ws.Range("D27").Value = "7.41"
$ws.Range("E27").Value = "  +4.98%  "
$ws.Range("D28").Value = "0.169"
$ws.Range("E28").Value = "  +5.72%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "26.07"
$ws.Range("E30").Value = "  +2.96%  "
$ws.Range("D31").Value = "0.109"
$ws.Range("E31").Value = "  +1.33%  "
$ws.Range("D32").Value = "10.37"
$ws.Range("E32").Value = "  +4.80%  "
$ws.Range("D33").Value = "51.34"
$ws.Range("E33").Value = "  +1.82%  "
$ws.Range("D34").Value = "34.44"
$ws.Range("E34").Value = "  +4.67%  "
$ws.Range("D35").Value = "2.06"
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("E36").Value = "  +3.33%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "3.29"
$ws.Range("E38").Value = "  +7.51%  "
$ws.Range("D39").Value = "17.15"
$ws.Range("E39").Value = "  +5.31%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "2.58"
$ws.Range("E40").Value = "  +7.02%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "1.85"
$ws.Range("E41").Value = "  +4.25%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.116"
$ws.Range("E42").Value = "  +1.23%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "127.05"
$ws.Range("E43").Value = "  +6.68%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "0.277"
$ws.Range("E44").Value = "  +8.14%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "3.82"
$ws.Range("E45").Value = "  +13.64%  "
$ws.Range("D46").Value = "21.43"
$ws.Range("E46").Value = "  +2.95%  "
$ws.Range("D47").Value = "2.05"
$ws.Range("E47").Value = "  +1.28%  "
$ws.Range("E48").Value = "  +1.04%  "
$ws.Range("D49").Value = "2.036.74"
$ws.Range("E49").Value = "  +2.42%  "
$ws.Range("D50").Value = "3.290.89"
$ws.Range("E50").Value = "  +2.42%  "
$ws.Range("D51").Value = "0.536"
$ws.Range("E51").Value = "  +20.75%  "
